$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Date of Birth" values for rows 3-26 in column L, continuing the
# sequential serial-date series started in L2 (33442 -> 1991-07-23),
# mirroring the same date format already applied to L2.
$srcCell = $ws.Range("L2")

for ($row = 3; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 12)
    $cell.Value = 33442 + ($row - 2)
    $cell.NumberFormat = $srcCell.NumberFormat
}

# Update the active selection to N18 (was L12).
$ws.Range("N18").Select()
